$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wheat's sowing time end date (row 10, column D = "End") updated
# from "31st December" to "30th November"
$ws.Range("D10").Value = "30th November"

# Reflect the resulting active cell/selection from the edit
$ws.Range("E10").Select()
